# Auto-generated edit script: applies "Add data for 2025-09-20" updates
# to the 2025 (column L) values across the Citywide Totals, By Neighborhood,
# and individual neighborhood worksheets. Only column L (year 2025) cells change.

$wb = $excel.ActiveWorkbook

# Sheet: Citywide Totals
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 4869  # was 4851
$ws.Range("L3").Value = 5229  # was 5205
$ws.Range("L4").Value = 1278  # was 1274
$ws.Range("L6").Value = 4424  # was 4412
$ws.Range("L7").Value = 16106  # was 16048

# Sheet: Austin
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 314  # was 311
$ws.Range("L3").Value = 363  # was 361
$ws.Range("L6").Value = 278  # was 277
$ws.Range("L7").Value = 1074  # was 1068

# Sheet: Garfield Park
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L3").Value = 255  # was 252
$ws.Range("L4").Value = 45  # was 44
$ws.Range("L7").Value = 738  # was 734

# Sheet: West Pullman
$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L3").Value = 69  # was 68
$ws.Range("L7").Value = 222  # was 221

# Sheet: Grand Crossing
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L3").Value = 205  # was 204
$ws.Range("L6").Value = 170  # was 169
$ws.Range("L7").Value = 608  # was 606

# Sheet: New City
$ws = $wb.Worksheets.Item("New City")
$ws.Range("L2").Value = 113  # was 112
$ws.Range("L7").Value = 314  # was 313

# Sheet: Woodlawn
$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L3").Value = 118  # was 116
$ws.Range("L7").Value = 276  # was 274

# Sheet: By Neighborhood
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L2").Value = 134  # was 132
$ws.Range("L8").Value = 1074  # was 1068
$ws.Range("L10").Value = 106  # was 105
$ws.Range("L11").Value = 261  # was 260
$ws.Range("L15").Value = 120  # was 118
$ws.Range("L20").Value = 403  # was 401
$ws.Range("L23").Value = 174  # was 173
$ws.Range("L25").Value = 95  # was 94
$ws.Range("L27").Value = 145  # was 143
$ws.Range("L29").Value = 880  # was 877
$ws.Range("L31").Value = 164  # was 163
$ws.Range("L33").Value = 738  # was 734
$ws.Range("L36").Value = 210  # was 209
$ws.Range("L37").Value = 608  # was 606
$ws.Range("L39").Value = 9  # was 10
$ws.Range("L42").Value = 525  # was 521
$ws.Range("L43").Value = 118  # was 117
$ws.Range("L44").Value = 112  # was 111
$ws.Range("L48").Value = 209  # was 208
$ws.Range("L50").Value = 80  # was 79
$ws.Range("L51").Value = 204  # was 203
$ws.Range("L52").Value = 325  # was 324
$ws.Range("L54").Value = 338  # was 336
$ws.Range("L55").Value = 152  # was 151
$ws.Range("L60").Value = 103  # was 102
$ws.Range("L63").Value = 43  # was 47
$ws.Range("L64").Value = 108  # was 109
$ws.Range("L65").Value = 314  # was 313
$ws.Range("L67").Value = 556  # was 553
$ws.Range("L75").Value = 58  # was 57
$ws.Range("L76").Value = 248  # was 247
$ws.Range("L77").Value = 106  # was 105
$ws.Range("L79").Value = 424  # was 422
$ws.Range("L84").Value = 157  # was 155
$ws.Range("L85").Value = 825  # was 824
$ws.Range("L86").Value = 116  # was 115
$ws.Range("L88").Value = 174  # was 173
$ws.Range("L90").Value = 162  # was 161
$ws.Range("L92").Value = 46  # was 45
$ws.Range("L95").Value = 222  # was 221
$ws.Range("L96").Value = 181  # was 180
$ws.Range("L98").Value = 88  # was 87
$ws.Range("L99").Value = 276  # was 274
$ws.Range("L101").Value = 16106  # was 16048

# Sheet: Gage Park
$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("L6").Value = 46  # was 45
$ws.Range("L7").Value = 164  # was 163

# Sheet: North Lawndale
$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L2").Value = 163  # was 162
$ws.Range("L3").Value = 214  # was 213
$ws.Range("L6").Value = 128  # was 127
$ws.Range("L7").Value = 556  # was 553

# Sheet: South Deering
$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L2").Value = 50  # was 49
$ws.Range("L6").Value = 46  # was 45
$ws.Range("L7").Value = 157  # was 155

# Sheet: Loop
$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L3").Value = 83  # was 82
$ws.Range("L4").Value = 29  # was 28
$ws.Range("L7").Value = 338  # was 336

# Sheet: Englewood
$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 266  # was 263
$ws.Range("L4").Value = 41  # was 40
$ws.Range("L6").Value = 228  # was 229
$ws.Range("L7").Value = 880  # was 877

# Sheet: Lake View
$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L4").Value = 42  # was 41
$ws.Range("L7").Value = 209  # was 208

# Sheet: Irving Park
$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("L3").Value = 32  # was 31
$ws.Range("L7").Value = 112  # was 111

# Sheet: River North
$ws = $wb.Worksheets.Item("River North")
$ws.Range("L2").Value = 50  # was 49
$ws.Range("L7").Value = 248  # was 247

# Sheet: Humboldt Park
$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L2").Value = 152  # was 151
$ws.Range("L3").Value = 175  # was 173
$ws.Range("L4").Value = 40  # was 39
$ws.Range("L7").Value = 525  # was 521

# Sheet: Avondale
$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("L3").Value = 29  # was 28
$ws.Range("L7").Value = 106  # was 105

# Sheet: Lower West Side
$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("L4").Value = 11  # was 10
$ws.Range("L7").Value = 152  # was 151

# Sheet: Douglas
$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L2").Value = 43  # was 42
$ws.Range("L7").Value = 174  # was 173

# Sheet: West Ridge
$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L2").Value = 56  # was 55
$ws.Range("L7").Value = 181  # was 180

# Sheet: Roseland
$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L3").Value = 151  # was 150
$ws.Range("L6").Value = 93  # was 92
$ws.Range("L7").Value = 424  # was 422

# Sheet: Near South Side
$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("L4").Value = 15  # was 16
$ws.Range("L7").Value = 108  # was 109

# Sheet: Chicago Lawn
$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L3").Value = 131  # was 129
$ws.Range("L7").Value = 403  # was 401

# Sheet: Grand Boulevard
$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L6").Value = 54  # was 53
$ws.Range("L7").Value = 210  # was 209

# Sheet: East Side
$ws = $wb.Worksheets.Item("East Side")
$ws.Range("L6").Value = 14  # was 13
$ws.Range("L7").Value = 95  # was 94

# Sheet: Brighton Park
$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("L2").Value = 42  # was 41
$ws.Range("L3").Value = 40  # was 39
$ws.Range("L7").Value = 120  # was 118

# Sheet: Wicker Park
$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("L2").Value = 22  # was 21
$ws.Range("L7").Value = 88  # was 87

# Sheet: Lincoln Square
$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("L6").Value = 20  # was 19
$ws.Range("L7").Value = 80  # was 79

# Sheet: Greektown
$ws = $wb.Worksheets.Item("Greektown")
$ws.Range("L5").Value = 4  # was 5
$ws.Range("L6").Value = 9  # was 10

# Sheet: Belmont Cragin
$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L3").Value = 79  # was 78
$ws.Range("L7").Value = 261  # was 260

# Sheet: Albany Park
$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("L3").Value = 44  # was 43
$ws.Range("L6").Value = 33  # was 32
$ws.Range("L7").Value = 134  # was 132

# Sheet: West Elsdon
$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("L6").Value = 18  # was 17
$ws.Range("L7").Value = 46  # was 45

# Sheet: United Center
$ws = $wb.Worksheets.Item("United Center")
$ws.Range("L2").Value = 50  # was 49
$ws.Range("L7").Value = 174  # was 173

# Sheet: Edgewater
$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("L4").Value = 20  # was 19
$ws.Range("L6").Value = 41  # was 40
$ws.Range("L7").Value = 145  # was 143

# Sheet: Streeterville
$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("L2").Value = 17  # was 16
$ws.Range("L7").Value = 116  # was 115

# Sheet: Pullman
$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("L2").Value = 26  # was 25
$ws.Range("L7").Value = 58  # was 57

# Sheet: Washington Heights
$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("L3").Value = 48  # was 47
$ws.Range("L7").Value = 162  # was 161

# Sheet: Little Italy, UIC
$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L6").Value = 47  # was 46
$ws.Range("L7").Value = 204  # was 203

# Sheet: Morgan Park
$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("L6").Value = 26  # was 25
$ws.Range("L7").Value = 103  # was 102

# Sheet: Hyde Park
$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("L3").Value = 37  # was 36
$ws.Range("L7").Value = 118  # was 117

# Sheet: South Shore
$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L3").Value = 335  # was 333
$ws.Range("L6").Value = 174  # was 175
$ws.Range("L7").Value = 825  # was 824

# Sheet: Riverdale
$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("L2").Value = 37  # was 36
$ws.Range("L7").Value = 106  # was 105

# Sheet: Little Village
$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L4").Value = 22  # was 21
$ws.Range("L7").Value = 325  # was 324
